$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws_ALC.Range("H15").Value = 1553.2778
$ws_ALC.Range("I15").Value = 1553.2778
$ws_ALC.Range("K15").Value = 4659.8334
$ws_ALC.Range("M15").Value = -4490.8334

# ALC row 132
$ws_ALC.Range("H132").Value = 3126.8635
$ws_ALC.Range("I132").Value = 2864.1018
$ws_ALC.Range("K132").Value = 8592.305399999999
$ws_ALC.Range("M132").Value = -6062.305399999999

# ALC row 138
$ws_ALC.Range("H138").Value = 2129.375
$ws_ALC.Range("I138").Value = 1805.96
$ws_ALC.Range("J138").Value = 2480.913
$ws_ALC.Range("K138").Value = 5417.88
$ws_ALC.Range("L138").Value = 7442.739
$ws_ALC.Range("M138").Value = -277.8800000000001
$ws_ALC.Range("N138").Value = -17722.739

# ALC row 141
$ws_ALC.Range("H141").Value = 1358.1666
$ws_ALC.Range("I141").Value = 1153
$ws_ALC.Range("K141").Value = 3459
$ws_ALC.Range("M141").Value = 1721

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws_ARM.Range("H2").Value = 1016.2727
$ws_ARM.Range("I2").Value = 618
$ws_ARM.Range("K2").Value = 618
$ws_ARM.Range("M2").Value = -505

# ARM row 45
$ws_ARM.Range("H45").Value = 2078.182
$ws_ARM.Range("I45").Value = 1936
$ws_ARM.Range("K45").Value = 1936
$ws_ARM.Range("M45").Value = -1559

# ARM row 56
$ws_ARM.Range("H56").Value = 24499.5
$ws_ARM.Range("J56").Value = 24499.5
$ws_ARM.Range("L56").Value = 24499.5
$ws_ARM.Range("N56").Value = -25983.5

# ARM row 61
$ws_ARM.Range("H61").Value = 1967.3889
$ws_ARM.Range("J61").Value = 970
$ws_ARM.Range("L61").Value = 970
$ws_ARM.Range("N61").Value = -1394

# ARM row 116
$ws_ARM.Range("H116").Value = 1016.2727
$ws_ARM.Range("I116").Value = 618
$ws_ARM.Range("K116").Value = 618
$ws_ARM.Range("M116").Value = 1676

# ARM row 122
$ws_ARM.Range("H122").Value = 2999.276
$ws_ARM.Range("J122").Value = 3999.889
$ws_ARM.Range("L122").Value = 11999.667
$ws_ARM.Range("N122").Value = -16899.667

# ARM row 136
$ws_ARM.Range("H136").Value = 1967.3889
$ws_ARM.Range("J136").Value = 970
$ws_ARM.Range("L136").Value = 2910
$ws_ARM.Range("N136").Value = -8010

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws_BSM.Range("H3").Value = 1016.2727
$ws_BSM.Range("I3").Value = 618
$ws_BSM.Range("K3").Value = 618
$ws_BSM.Range("M3").Value = -504

# BSM row 5
$ws_BSM.Range("H5").Value = 40372
$ws_BSM.Range("J5").Value = 40372
$ws_BSM.Range("L5").Value = 40372
$ws_BSM.Range("N5").Value = -40598

# BSM row 22
$ws_BSM.Range("H22").Value = 125
$ws_BSM.Range("I22").Value = 125
$ws_BSM.Range("J22").Value = 0
$ws_BSM.Range("K22").Value = 125
$ws_BSM.Range("L22").Value = 0
$ws_BSM.Range("M22").Value = 48
$ws_BSM.Range("N22").ClearContents()

# BSM row 107
$ws_BSM.Range("H107").Value = 1842.9
$ws_BSM.Range("I107").Value = 1546.2307
$ws_BSM.Range("J107").Value = 2393.8572
$ws_BSM.Range("K107").Value = 1546.2307
$ws_BSM.Range("L107").Value = 2393.8572
$ws_BSM.Range("M107").Value = 373.7692999999999
$ws_BSM.Range("N107").Value = -6233.8572

# BSM row 134
$ws_BSM.Range("H134").Value = 1212
$ws_BSM.Range("I134").Value = 1119.7222
$ws_BSM.Range("K134").Value = 3359.1666
$ws_BSM.Range("M134").Value = -824.1665999999996

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws_CRP.Range("H16").Value = 2227.6
$ws_CRP.Range("I16").Value = 2207.875
$ws_CRP.Range("J16").Value = 2306.5
$ws_CRP.Range("K16").Value = 2207.875
$ws_CRP.Range("L16").Value = 2306.5
$ws_CRP.Range("M16").Value = -1920.875
$ws_CRP.Range("N16").Value = -2880.5

# CRP row 45
$ws_CRP.Range("H45").Value = 0
$ws_CRP.Range("J45").Value = 0
$ws_CRP.Range("L45").Value = 0
$ws_CRP.Range("N45").ClearContents()

# CRP row 58
$ws_CRP.Range("H58").Value = 2525.923
$ws_CRP.Range("J58").Value = 2731
$ws_CRP.Range("L58").Value = 2731
$ws_CRP.Range("N58").Value = -3137

# CRP row 105
$ws_CRP.Range("H105").Value = 3659.875
$ws_CRP.Range("I105").Value = 2855.8
$ws_CRP.Range("J105").Value = 5000
$ws_CRP.Range("K105").Value = 2855.8
$ws_CRP.Range("L105").Value = 5000
$ws_CRP.Range("M105").Value = -1108.8
$ws_CRP.Range("N105").Value = -8494

# CRP row 113
$ws_CRP.Range("H113").Value = 2227.6
$ws_CRP.Range("I113").Value = 2207.875
$ws_CRP.Range("J113").Value = 2306.5
$ws_CRP.Range("K113").Value = 2207.875
$ws_CRP.Range("L113").Value = 2306.5
$ws_CRP.Range("M113").Value = -37.875
$ws_CRP.Range("N113").Value = -6646.5

# CRP row 122
$ws_CRP.Range("H122").Value = 2638
$ws_CRP.Range("I122").Value = 2577.25
$ws_CRP.Range("K122").Value = 7731.75
$ws_CRP.Range("M122").Value = -5281.75

# CRP row 132
$ws_CRP.Range("H132").Value = 1871.8077
$ws_CRP.Range("I132").Value = 1310.2858
$ws_CRP.Range("J132").Value = 4230.2
$ws_CRP.Range("K132").Value = 3930.8574
$ws_CRP.Range("L132").Value = 12690.6
$ws_CRP.Range("M132").Value = -1400.8574
$ws_CRP.Range("N132").Value = -17750.6

# CRP row 134
$ws_CRP.Range("H134").Value = 2173.9375
$ws_CRP.Range("I134").Value = 1914.5714
$ws_CRP.Range("K134").Value = 5743.7142
$ws_CRP.Range("M134").Value = -3208.7142

# CRP row 136
$ws_CRP.Range("H136").Value = 2525.923
$ws_CRP.Range("J136").Value = 2731
$ws_CRP.Range("L136").Value = 8193
$ws_CRP.Range("N136").Value = -13293

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws_CUL.Range("H4").Value = 15705650
$ws_CUL.Range("I4").Value = 18003306
$ws_CUL.Range("K4").Value = 54009918
$ws_CUL.Range("M4").Value = -54009806

# CUL row 38
$ws_CUL.Range("H38").Value = 96.40909000000001
$ws_CUL.Range("I38").Value = 73.22221999999999
$ws_CUL.Range("K38").Value = 219.66666
$ws_CUL.Range("M38").Value = 127.33334

# CUL row 68
$ws_CUL.Range("H68").Value = 1498.2
$ws_CUL.Range("J68").Value = 1935.4
$ws_CUL.Range("L68").Value = 5806.200000000001
$ws_CUL.Range("N68").Value = -7428.200000000001

# CUL row 71
$ws_CUL.Range("H71").Value = 1498.2
$ws_CUL.Range("J71").Value = 1935.4
$ws_CUL.Range("L71").Value = 17418.6
$ws_CUL.Range("N71").Value = -25530.6

# CUL row 107
$ws_CUL.Range("H107").Value = 947.8570999999999
$ws_CUL.Range("J107").Value = 827.5
$ws_CUL.Range("L107").Value = 2482.5
$ws_CUL.Range("N107").Value = -6322.5

# CUL row 132
$ws_CUL.Range("H132").Value = 1761.5454
$ws_CUL.Range("J132").Value = 1673.75
$ws_CUL.Range("L132").Value = 15063.75
$ws_CUL.Range("N132").Value = -20123.75

# CUL row 140
$ws_CUL.Range("H140").Value = 2702.5
$ws_CUL.Range("I140").Value = 1132
$ws_CUL.Range("K140").Value = 3396
$ws_CUL.Range("M140").Value = 1784

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws_GSM.Range("H80").Value = 4082.2942
$ws_GSM.Range("I80").Value = 4612.222
$ws_GSM.Range("J80").Value = 3486.125
$ws_GSM.Range("K80").Value = 4612.222
$ws_GSM.Range("L80").Value = 3486.125
$ws_GSM.Range("M80").Value = -3614.222
$ws_GSM.Range("N80").Value = -5482.125

# GSM row 83
$ws_GSM.Range("H83").Value = 4082.2942
$ws_GSM.Range("I83").Value = 4612.222
$ws_GSM.Range("J83").Value = 3486.125
$ws_GSM.Range("K83").Value = 23061.11
$ws_GSM.Range("L83").Value = 17430.625
$ws_GSM.Range("M83").Value = -18069.11
$ws_GSM.Range("N83").Value = -27414.625

# GSM row 102
$ws_GSM.Range("H102").Value = 2315.75
$ws_GSM.Range("I102").Value = 1921.25
$ws_GSM.Range("J102").Value = 3499.25
$ws_GSM.Range("K102").Value = 1921.25
$ws_GSM.Range("L102").Value = 3499.25
$ws_GSM.Range("M102").Value = -299.25
$ws_GSM.Range("N102").Value = -6743.25

# GSM row 126
$ws_GSM.Range("H126").Value = 6815.391
$ws_GSM.Range("I126").Value = 8788.532999999999
$ws_GSM.Range("J126").Value = 3115.75
$ws_GSM.Range("K126").Value = 26365.599
$ws_GSM.Range("L126").Value = 9347.25
$ws_GSM.Range("M126").Value = -23895.599
$ws_GSM.Range("N126").Value = -14287.25

# GSM row 132
$ws_GSM.Range("H132").Value = 5532.4707
$ws_GSM.Range("I132").Value = 5186.4116
$ws_GSM.Range("J132").Value = 5878.5293
$ws_GSM.Range("K132").Value = 15559.2348
$ws_GSM.Range("L132").Value = 17635.5879
$ws_GSM.Range("M132").Value = -13029.2348
$ws_GSM.Range("N132").Value = -22695.5879

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 3281.9285
$ws_LTW.Range("I7").Value = 2594.7
$ws_LTW.Range("J7").Value = 5000
$ws_LTW.Range("K7").Value = 2594.7
$ws_LTW.Range("L7").Value = 5000
$ws_LTW.Range("M7").Value = -2482.7
$ws_LTW.Range("N7").Value = -5224

# LTW row 22
$ws_LTW.Range("H22").Value = 1702633.5
$ws_LTW.Range("J22").Value = 4042346
$ws_LTW.Range("L22").Value = 4042346
$ws_LTW.Range("N22").Value = -4042936

# LTW row 25
$ws_LTW.Range("H25").Value = 1025000
$ws_LTW.Range("I25").Value = 1025000
$ws_LTW.Range("K25").Value = 1025000
$ws_LTW.Range("M25").Value = -1024770

# LTW row 27
$ws_LTW.Range("H27").Value = 1702633.5
$ws_LTW.Range("J27").Value = 4042346
$ws_LTW.Range("L27").Value = 4042346
$ws_LTW.Range("N27").Value = -4042560

# LTW row 40
$ws_LTW.Range("H40").Value = 3171.2856
$ws_LTW.Range("J40").Value = 4499.75
$ws_LTW.Range("L40").Value = 4499.75
$ws_LTW.Range("N40").Value = -4771.75

# LTW row 126
$ws_LTW.Range("H126").Value = 3281.9285
$ws_LTW.Range("I126").Value = 2594.7
$ws_LTW.Range("J126").Value = 5000
$ws_LTW.Range("K126").Value = 7784.099999999999
$ws_LTW.Range("L126").Value = 15000
$ws_LTW.Range("M126").Value = -5314.099999999999
$ws_LTW.Range("N126").Value = -19940

# LTW row 132
$ws_LTW.Range("H132").Value = 6550.913
$ws_LTW.Range("I132").Value = 3704.4167
$ws_LTW.Range("J132").Value = 9656.182000000001
$ws_LTW.Range("K132").Value = 11113.2501
$ws_LTW.Range("L132").Value = 28968.546
$ws_LTW.Range("M132").Value = -8583.250100000001
$ws_LTW.Range("N132").Value = -34028.546

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 18
$ws_WVR.Range("H18").Value = 2699.8572
$ws_WVR.Range("I18").Value = 1000
$ws_WVR.Range("J18").Value = 12899
$ws_WVR.Range("K18").Value = 1000
$ws_WVR.Range("L18").Value = 12899
$ws_WVR.Range("M18").Value = -827
$ws_WVR.Range("N18").Value = -13245

# WVR row 70
$ws_WVR.Range("H70").Value = 35766.668
$ws_WVR.Range("I70").Value = 0
$ws_WVR.Range("J70").Value = 35766.668
$ws_WVR.Range("K70").Value = 0
$ws_WVR.Range("L70").Value = 35766.668
$ws_WVR.Range("M70").ClearContents()
$ws_WVR.Range("N70").Value = -36396.668

# WVR row 73
$ws_WVR.Range("H73").Value = 35766.668
$ws_WVR.Range("I73").Value = 0
$ws_WVR.Range("J73").Value = 35766.668
$ws_WVR.Range("K73").Value = 0
$ws_WVR.Range("L73").Value = 35766.668
$ws_WVR.Range("M73").ClearContents()
$ws_WVR.Range("N73").Value = -37950.668

# WVR row 122
$ws_WVR.Range("H122").Value = 1932.9259
$ws_WVR.Range("I122").Value = 1768.0526
$ws_WVR.Range("J122").Value = 2324.5
$ws_WVR.Range("K122").Value = 5304.1578
$ws_WVR.Range("L122").Value = 6973.5
$ws_WVR.Range("M122").Value = -2854.1578
$ws_WVR.Range("N122").Value = -11873.5

# WVR row 136
$ws_WVR.Range("H136").Value = 2941.7334
$ws_WVR.Range("I136").Value = 1420.5
$ws_WVR.Range("K136").Value = 4261.5
$ws_WVR.Range("M136").Value = -1711.5
